$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.100.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.925.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.79%  '

$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.21'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5081'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.79%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4033'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08374'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.119'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.34%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.23'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.04%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.07%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.411'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.923.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.273'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.27%  '

$ws.Range("E16").Value = '  -0.02%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.94'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001098'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06498'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.37%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.51'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.953'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.103.28'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.01%  '

$ws.Range("E25").Value = '  +0.63%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.143.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.90'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.07%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.266'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.82%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.99'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.30%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.133'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.45%  '

$ws.Range("E32").Value = '  +1.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.965'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.45%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.787'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.71%  '

$ws.Range("E35").Value = '  +1.07%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.265'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.61%  '

$ws.Range("B37").Value = 'InternetComputer(DFINITY)'
$ws.Range("C37").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.308'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.35%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06457'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2152'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '

$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.74'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.53%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6484'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.83%  '

$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.691'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.21%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.220'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.32'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.25%  '

$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6066'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.165'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.59%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.631'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.45'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.211'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '

$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '78.07'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.83%  '

$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.130'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.77%  '
